# Update the results sheet:
#  - rename H1 "Experiment_Time" -> "Avg_Experiment_Time"
#  - insert 6 new "Std_*" columns (I..N) before the existing "Obs_Prob" column,
#    which moves from I to O
#  - the existing "Std_Total_Rounds" column (old J) is recomputed and moved to
#    new column I
#  - refresh every data row (2-13) with the recomputed Avg_* / Std_* values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row -----------------------------------------------------------
$ws.Cells.Item(1, 8).Value  = "Avg_Experiment_Time"
$ws.Cells.Item(1, 9).Value  = "Std_Total_Rounds"
$ws.Cells.Item(1, 10).Value = "Std_Expl_Cost"
$ws.Cells.Item(1, 11).Value = "Std_Expl_Eff"
$ws.Cells.Item(1, 12).Value = "Std_Round_Time"
$ws.Cells.Item(1, 13).Value = "Std_Agent_Step_Time"
$ws.Cells.Item(1, 14).Value = "Std_Experiment_Time"
$ws.Cells.Item(1, 15).Value = "Obs_Prob"

# the freshly-created header cells (K1:O1) need the same bold/centered/
# bordered look as the rest of row 1 (style index 1) - copy it over from an
# existing header cell
$ws.Range("J1").Copy()
$ws.Range("K1:O1").PasteSpecial(-4122)

# ---- data rows -------------------------------------------------------------
# columns: A=#_Agents B=Coverage C=Avg_Total_Rounds D=Avg_Expl_Cost
#          E=Avg_Expl_Eff F=Avg_Round_Time G=Avg_Agent_Step_Time
#          H=Avg_Experiment_Time I=Std_Total_Rounds J=Std_Expl_Cost
#          K=Std_Expl_Eff L=Std_Round_Time M=Std_Agent_Step_Time
#          N=Std_Experiment_Time O=Obs_Prob

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 221.864
$ws.Cells.Item(2, 4).Value = 221.864
$ws.Cells.Item(2, 5).Value = 0.8634281
$ws.Cells.Item(2, 6).Value = 0.02250744
$ws.Cells.Item(2, 7).Value = 0.02250744
$ws.Cells.Item(2, 8).Value = 4.961661459999999
$ws.Cells.Item(2, 9).Value = 83.83915261393119
$ws.Cells.Item(2, 10).Value = 83.83915261393119
$ws.Cells.Item(2, 11).Value = 0.285033948210129
$ws.Cells.Item(2, 12).Value = 0.001083746862012455
$ws.Cells.Item(2, 13).Value = 0.001083746862012455
$ws.Cells.Item(2, 14).Value = 1.796676182025739
$ws.Cells.Item(2, 15).Value = 0.15

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 308.98
$ws.Cells.Item(3, 4).Value = 308.98
$ws.Cells.Item(3, 5).Value = 0.59875896
$ws.Cells.Item(3, 6).Value = 0.008194459999999999
$ws.Cells.Item(3, 7).Value = 0.008194459999999999
$ws.Cells.Item(3, 8).Value = 2.48084538
$ws.Cells.Item(3, 9).Value = 93.76809965561877
$ws.Cells.Item(3, 10).Value = 93.76809965561877
$ws.Cells.Item(3, 11).Value = 0.1873628341656144
$ws.Cells.Item(3, 12).Value = 0.001949715437949332
$ws.Cells.Item(3, 13).Value = 0.001949715437949332
$ws.Cells.Item(3, 14).Value = 0.8065002478611248
$ws.Cells.Item(3, 15).Value = 0.85

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 119.09
$ws.Cells.Item(4, 4).Value = 238.18
$ws.Cells.Item(4, 5).Value = 0.78331358
$ws.Cells.Item(4, 6).Value = 0.0442501
$ws.Cells.Item(4, 7).Value = 0.02212494
$ws.Cells.Item(4, 8).Value = 2.62119988
$ws.Cells.Item(4, 9).Value = 39.46020927760603
$ws.Cells.Item(4, 10).Value = 78.92041855521205
$ws.Cells.Item(4, 11).Value = 0.2379316683238411
$ws.Cells.Item(4, 12).Value = 0.0026279551912036
$ws.Cells.Item(4, 13).Value = 0.001314125278234076
$ws.Cells.Item(4, 14).Value = 0.8325841288956711
$ws.Cells.Item(4, 15).Value = 0.15

$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = 161.502
$ws.Cells.Item(5, 4).Value = 322.986
$ws.Cells.Item(5, 5).Value = 0.5692335799999999
$ws.Cells.Item(5, 6).Value = 0.01418536
$ws.Cells.Item(5, 7).Value = 0.007092520000000001
$ws.Cells.Item(5, 8).Value = 1.11636568
$ws.Cells.Item(5, 9).Value = 47.08302899387928
$ws.Cells.Item(5, 10).Value = 94.16030017080719
$ws.Cells.Item(5, 11).Value = 0.1696750924058965
$ws.Cells.Item(5, 12).Value = 0.003521931298384607
$ws.Cells.Item(5, 13).Value = 0.00176086972089019
$ws.Cells.Item(5, 14).Value = 0.3387495214525125
$ws.Cells.Item(5, 15).Value = 0.85

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = 62.008
$ws.Cells.Item(6, 4).Value = 248.024
$ws.Cells.Item(6, 5).Value = 0.7404575800000001
$ws.Cells.Item(6, 6).Value = 0.08357519999999999
$ws.Cells.Item(6, 7).Value = 0.02089374
$ws.Cells.Item(6, 8).Value = 1.29277098
$ws.Cells.Item(6, 9).Value = 17.92558476531119
$ws.Cells.Item(6, 10).Value = 71.70362783352721
$ws.Cells.Item(6, 11).Value = 0.2431236364674606
$ws.Cells.Item(6, 12).Value = 0.003689096885643331
$ws.Cells.Item(6, 13).Value = 0.0009222816376999364
$ws.Cells.Item(6, 14).Value = 0.3663289897246257
$ws.Cells.Item(6, 15).Value = 0.15

$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = 84.136
$ws.Cells.Item(7, 4).Value = 336.366
$ws.Cells.Item(7, 5).Value = 0.5452307800000001
$ws.Cells.Item(7, 6).Value = 0.02249988
$ws.Cells.Item(7, 7).Value = 0.00562472
$ws.Cells.Item(7, 8).Value = 0.46655192
$ws.Cells.Item(7, 9).Value = 23.97638918763528
$ws.Cells.Item(7, 10).Value = 95.87289119731952
$ws.Cells.Item(7, 11).Value = 0.1632305901340267
$ws.Cells.Item(7, 12).Value = 0.006792230707327537
$ws.Cells.Item(7, 13).Value = 0.001698099999651857
$ws.Cells.Item(7, 14).Value = 0.1780225400589816
$ws.Cells.Item(7, 15).Value = 0.85

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = 44.696
$ws.Cells.Item(8, 4).Value = 268.168
$ws.Cells.Item(8, 5).Value = 0.68076348
$ws.Cells.Item(8, 6).Value = 0.12028286
$ws.Cells.Item(8, 7).Value = 0.02004722
$ws.Cells.Item(8, 8).Value = 0.8954996200000001
$ws.Cells.Item(8, 9).Value = 11.92089760039899
$ws.Cells.Item(8, 10).Value = 71.51884833618806
$ws.Cells.Item(8, 11).Value = 0.2354337352191654
$ws.Cells.Item(8, 12).Value = 0.007124428467067708
$ws.Cells.Item(8, 13).Value = 0.00118737128681196
$ws.Cells.Item(8, 14).Value = 0.2418541251095805
$ws.Cells.Item(8, 15).Value = 0.15

$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 56.762
$ws.Cells.Item(9, 4).Value = 340.226
$ws.Cells.Item(9, 5).Value = 0.5448089600000001
$ws.Cells.Item(9, 6).Value = 0.0258413
$ws.Cells.Item(9, 7).Value = 0.004306940000000001
$ws.Cells.Item(9, 8).Value = 0.24364626
$ws.Cells.Item(9, 9).Value = 15.81117777666707
$ws.Cells.Item(9, 10).Value = 94.78097800499405
$ws.Cells.Item(9, 11).Value = 0.1944558111631242
$ws.Cells.Item(9, 12).Value = 0.008418955120479101
$ws.Cells.Item(9, 13).Value = 0.001403227469465348
$ws.Cells.Item(9, 14).Value = 0.1022634210398711
$ws.Cells.Item(9, 15).Value = 0.85

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = 34.936
$ws.Cells.Item(10, 4).Value = 279.444
$ws.Cells.Item(10, 5).Value = 0.65603296
$ws.Cells.Item(10, 6).Value = 0.15256754
$ws.Cells.Item(10, 7).Value = 0.01907088
$ws.Cells.Item(10, 8).Value = 0.66530704
$ws.Cells.Item(10, 9).Value = 8.702607353517623
$ws.Cells.Item(10, 10).Value = 69.62486702476662
$ws.Cells.Item(10, 11).Value = 0.25934858865572
$ws.Cells.Item(10, 12).Value = 0.01053728831523922
$ws.Cells.Item(10, 13).Value = 0.001317214688146673
$ws.Cells.Item(10, 14).Value = 0.1674058303647714
$ws.Cells.Item(10, 15).Value = 0.15

$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(11, 3).Value = 44.392
$ws.Cells.Item(11, 4).Value = 354.538
$ws.Cells.Item(11, 5).Value = 0.5381165600000001
$ws.Cells.Item(11, 6).Value = 0.0274612
$ws.Cells.Item(11, 7).Value = 0.00343254
$ws.Cells.Item(11, 8).Value = 0.15098592
$ws.Cells.Item(11, 9).Value = 14.28518354582383
$ws.Cells.Item(11, 10).Value = 114.1176332757076
$ws.Cells.Item(11, 11).Value = 0.21722963059164
$ws.Cells.Item(11, 12).Value = 0.01119122703214308
$ws.Cells.Item(11, 13).Value = 0.00139881647239618
$ws.Cells.Item(11, 14).Value = 0.0775893894735056
$ws.Cells.Item(11, 15).Value = 0.85

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = 29.254
$ws.Cells.Item(12, 4).Value = 292.426
$ws.Cells.Item(12, 5).Value = 0.6599942800000002
$ws.Cells.Item(12, 6).Value = 0.18216406
$ws.Cells.Item(12, 7).Value = 0.0182164
$ws.Cells.Item(12, 8).Value = 0.5292801400000001
$ws.Cells.Item(12, 9).Value = 7.880716680022291
$ws.Cells.Item(12, 10).Value = 78.7595496261541
$ws.Cells.Item(12, 11).Value = 0.4063990917667746
$ws.Cells.Item(12, 12).Value = 0.02434057779863787
$ws.Cells.Item(12, 13).Value = 0.002434139665636228
$ws.Cells.Item(12, 14).Value = 0.1426137234027146
$ws.Cells.Item(12, 15).Value = 0.15

$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = 34.618
$ws.Cells.Item(13, 4).Value = 345.352
$ws.Cells.Item(13, 5).Value = 0.5832312000000001
$ws.Cells.Item(13, 6).Value = 0.0276596
$ws.Cells.Item(13, 7).Value = 0.0027658
$ws.Cells.Item(13, 8).Value = 0.09507578000000001
$ws.Cells.Item(13, 9).Value = 12.60610581454936
$ws.Cells.Item(13, 10).Value = 125.8041074165934
$ws.Cells.Item(13, 11).Value = 0.302567425215874
$ws.Cells.Item(13, 12).Value = 0.01193070595266928
$ws.Cells.Item(13, 13).Value = 0.001193086447220436
$ws.Cells.Item(13, 14).Value = 0.05407955028946359
$ws.Cells.Item(13, 15).Value = 0.85
